# Insert a new data row at row 394 (pushes existing rows 394:491 down to 395:492)
# and populate it with the new price-record values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(394).Insert()

$ws.Range("A394").Value = 10
$ws.Range("B394").Value = "Vega Modelo de Temuco"
$ws.Range("C394").Value = "La Araucanía"
$ws.Range("D394").Value = 44889
$ws.Range("E394").Value = 9
$ws.Range("F394").Value = 100114014
$ws.Range("G394").Value = "Betarraga"
$ws.Range("H394").Value = "Sin especificar"
$ws.Range("I394").Value = "Primera"
$ws.Range("J394").Value = 55
$ws.Range("K394").Value = 11000
$ws.Range("L394").Value = 11000
$ws.Range("M394").Value = 11000
$ws.Range("N394").Value = "`$/docena de paquetes"
$ws.Range("O394").Value = "Provincia de Cautín"
$ws.Range("P394").Value = 917
$ws.Range("Q394").Value = 12
$ws.Range("R394").Value = "Hortaliza"
